$d = $word.ActiveDocument

$found = $d.Content.Find.Execute(
    "[ ] Affiliate disclosure (prominent)^p",
    $true, $false, $false, $false, $false,
    $true, 1, $false, "", 2)

Write-Output "Find executed, found=$found"
